$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B9: store phone number as a real number instead of text
$ws.Range("B9").Value = 53023029302

# Add new row 10 with Pedro Pablo's data
$ws.Range("A10").Value = "Pedro Pablo"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "38434734"
$ws.Range("C10").Value = "pp@gmail.com"
